$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")

$ws_ALC.Range("H53").Value = 1760.3636
$ws_ALC.Range("I53").Value = 3112
$ws_ALC.Range("J53").Value = 138.4
$ws_ALC.Range("K53").Value = 3112
$ws_ALC.Range("L53").Value = 138.4
$ws_ALC.Range("M53").Value = -2475
$ws_ALC.Range("N53").Value = -1412.4
$ws_ALC.Range("H70").Value = 1148.5714
$ws_ALC.Range("I70").Value = 887.8333
$ws_ALC.Range("K70").Value = 2663.4999
$ws_ALC.Range("M70").Value = -2393.4999
$ws_ALC.Range("H73").Value = 1148.5714
$ws_ALC.Range("I73").Value = 887.8333
$ws_ALC.Range("K73").Value = 2663.4999
$ws_ALC.Range("M73").Value = -1727.4999
$ws_ALC.Range("H82").Value = 964.13336
$ws_ALC.Range("I82").Value = 964.13336
$ws_ALC.Range("K82").Value = 2892.40008
$ws_ALC.Range("M82").Value = -2486.40008
$ws_ALC.Range("H85").Value = 964.13336
$ws_ALC.Range("I85").Value = 964.13336
$ws_ALC.Range("K85").Value = 2892.40008
$ws_ALC.Range("M85").Value = -1488.40008
$ws_ALC.Range("H86").Value = 7661.6665
$ws_ALC.Range("I86").Value = 7614.2
$ws_ALC.Range("J86").Value = 7721
$ws_ALC.Range("K86").Value = 7614.2
$ws_ALC.Range("L86").Value = 7721
$ws_ALC.Range("M86").Value = -6491.2
$ws_ALC.Range("N86").Value = -9967
$ws_ALC.Range("H88").Value = 3372.8948
$ws_ALC.Range("J88").Value = 3098.6155
$ws_ALC.Range("L88").Value = 3098.6155
$ws_ALC.Range("N88").Value = -3910.6155
$ws_ALC.Range("H89").Value = 7661.6665
$ws_ALC.Range("I89").Value = 7614.2
$ws_ALC.Range("J89").Value = 7721
$ws_ALC.Range("K89").Value = 38071
$ws_ALC.Range("L89").Value = 38605
$ws_ALC.Range("M89").Value = -32455
$ws_ALC.Range("N89").Value = -49837
$ws_ALC.Range("H91").Value = 3372.8948
$ws_ALC.Range("J91").Value = 3098.6155
$ws_ALC.Range("L91").Value = 3098.6155
$ws_ALC.Range("N91").Value = -5906.6155
$ws_ALC.Range("H100").Value = 3580.5
$ws_ALC.Range("I100").Value = 2800
$ws_ALC.Range("J100").Value = 3840.6667
$ws_ALC.Range("K100").Value = 2800
$ws_ALC.Range("L100").Value = 3840.6667
$ws_ALC.Range("M100").Value = -2259
$ws_ALC.Range("N100").Value = -4922.6667
$ws_ALC.Range("H116").Value = 7467.4443
$ws_ALC.Range("I116").Value = 2240
$ws_ALC.Range("J116").Value = 14001.75
$ws_ALC.Range("K116").Value = 2240
$ws_ALC.Range("L116").Value = 14001.75
$ws_ALC.Range("M116").Value = 1202
$ws_ALC.Range("N116").Value = -20885.75
$ws_ALC.Range("H135").Value = 29412758
$ws_ALC.Range("I135").Value = 957.76666
$ws_ALC.Range("J135").Value = 250001260
$ws_ALC.Range("K135").Value = 8619.899939999999
$ws_ALC.Range("L135").Value = 2250011340
$ws_ALC.Range("M135").Value = -6084.899939999999
$ws_ALC.Range("N135").Value = -2250016410
$ws_ARM = $wb.Worksheets.Item("ARM")

$ws_ARM.Range("H88").Value = 2116.2354
$ws_ARM.Range("J88").Value = 2090.4
$ws_ARM.Range("L88").Value = 2090.4
$ws_ARM.Range("N88").Value = -2902.4
$ws_ARM.Range("H91").Value = 2116.2354
$ws_ARM.Range("J91").Value = 2090.4
$ws_ARM.Range("L91").Value = 2090.4
$ws_ARM.Range("N91").Value = -4898.4
$ws_BSM = $wb.Worksheets.Item("BSM")

$ws_BSM.Range("H86").Value = 1857.8182
$ws_BSM.Range("I86").Value = 1755.5385
$ws_BSM.Range("J86").Value = 2005.5555
$ws_BSM.Range("K86").Value = 1755.5385
$ws_BSM.Range("L86").Value = 2005.5555
$ws_BSM.Range("M86").Value = -632.5385000000001
$ws_BSM.Range("N86").Value = -4251.5555
$ws_BSM.Range("H89").Value = 1857.8182
$ws_BSM.Range("I89").Value = 1755.5385
$ws_BSM.Range("J89").Value = 2005.5555
$ws_BSM.Range("K89").Value = 8777.692500000001
$ws_BSM.Range("L89").Value = 10027.7775
$ws_BSM.Range("M89").Value = -3161.692500000001
$ws_BSM.Range("N89").Value = -21259.7775
$ws_BSM.Range("H105").Value = 2344.9387
$ws_BSM.Range("I105").Value = 2400
$ws_BSM.Range("J105").Value = 2315.6875
$ws_BSM.Range("K105").Value = 2400
$ws_BSM.Range("L105").Value = 2315.6875
$ws_BSM.Range("M105").Value = -653
$ws_BSM.Range("N105").Value = -5809.6875
$ws_CRP = $wb.Worksheets.Item("CRP")

$ws_CRP.Range("H4").Value = 3379.8
$ws_CRP.Range("I4").Value = 0
$ws_CRP.Range("J4").Value = 3379.8
$ws_CRP.Range("K4").Value = 0
$ws_CRP.Range("M4").ClearContents()
$ws_CRP.Range("N4").Value = -3603.8
$ws_CRP.Range("H62").Value = 2301.375
$ws_CRP.Range("I62").Value = 2272.8333
$ws_CRP.Range("J62").Value = 2387
$ws_CRP.Range("K62").Value = 2272.8333
$ws_CRP.Range("L62").Value = 2387
$ws_CRP.Range("M62").Value = -1648.8333
$ws_CRP.Range("N62").Value = -3635
$ws_CRP.Range("H65").Value = 2301.375
$ws_CRP.Range("I65").Value = 2272.8333
$ws_CRP.Range("J65").Value = 2387
$ws_CRP.Range("K65").Value = 11364.1665
$ws_CRP.Range("L65").Value = 11935
$ws_CRP.Range("M65").Value = -8244.166499999999
$ws_CRP.Range("N65").Value = -18175
$ws_CRP.Range("H105").Value = 961.5
$ws_CRP.Range("I105").Value = 967.2727
$ws_CRP.Range("J105").Value = 940.3333
$ws_CRP.Range("K105").Value = 967.2727
$ws_CRP.Range("L105").Value = 940.3333
$ws_CRP.Range("M105").Value = 779.7273
$ws_CRP.Range("N105").Value = -4434.3333
$ws_CUL = $wb.Worksheets.Item("CUL")

$ws_CUL.Range("H80").Value = 662
$ws_CUL.Range("J80").Value = 662
$ws_CUL.Range("L80").Value = 1986
$ws_CUL.Range("N80").Value = -3858
$ws_CUL.Range("H83").Value = 662
$ws_CUL.Range("J83").Value = 662
$ws_CUL.Range("L83").Value = 5958
$ws_CUL.Range("N83").Value = -15318
$ws_CUL.Range("H99").Value = 1124.3
$ws_CUL.Range("I99").Value = 1124.3
$ws_CUL.Range("K99").Value = 3372.9
$ws_CUL.Range("M99").Value = -1126.9
$ws_CUL.Range("H113").Value = 1443501.4
$ws_CUL.Range("I113").Value = 3367523.5
$ws_CUL.Range("J113").Value = 484.91666
$ws_CUL.Range("K113").Value = 10102570.5
$ws_CUL.Range("L113").Value = 1454.74998
$ws_CUL.Range("M113").Value = -10100400.5
$ws_CUL.Range("N113").Value = -5794.749980000001
$ws_CUL.Range("H131").Value = 916.45
$ws_CUL.Range("J131").Value = 916.45
$ws_CUL.Range("L131").Value = 2749.35
$ws_CUL.Range("N131").Value = -12829.35
$ws_GSM = $wb.Worksheets.Item("GSM")

$ws_GSM.Range("H5").Value = 7898.6
$ws_GSM.Range("J5").Value = 10332
$ws_GSM.Range("L5").Value = 10332
$ws_GSM.Range("N5").Value = -10556
$ws_GSM.Range("H80").Value = 2617.3
$ws_GSM.Range("I80").Value = 2342
$ws_GSM.Range("J80").Value = 3128.5715
$ws_GSM.Range("K80").Value = 2342
$ws_GSM.Range("L80").Value = 3128.5715
$ws_GSM.Range("M80").Value = -1344
$ws_GSM.Range("N80").Value = -5124.5715
$ws_GSM.Range("H83").Value = 2617.3
$ws_GSM.Range("I83").Value = 2342
$ws_GSM.Range("J83").Value = 3128.5715
$ws_GSM.Range("K83").Value = 11710
$ws_GSM.Range("L83").Value = 15642.8575
$ws_GSM.Range("M83").Value = -6718
$ws_GSM.Range("N83").Value = -25626.8575
$ws_GSM.Range("H122").Value = 32348.656
$ws_GSM.Range("I122").Value = 38698.855
$ws_GSM.Range("K122").Value = 116096.565
$ws_GSM.Range("M122").Value = -113646.565
$ws_LTW = $wb.Worksheets.Item("LTW")

$ws_LTW.Range("H7").Value = 2624.75
$ws_LTW.Range("I7").Value = 2500
$ws_LTW.Range("J7").Value = 2999
$ws_LTW.Range("K7").Value = 2500
$ws_LTW.Range("L7").Value = 2999
$ws_LTW.Range("M7").Value = -2388
$ws_LTW.Range("N7").Value = -3223
$ws_LTW.Range("H68").Value = 8158.684
$ws_LTW.Range("I68").Value = 22078
$ws_LTW.Range("J68").Value = 3187.5
$ws_LTW.Range("K68").Value = 22078
$ws_LTW.Range("L68").Value = 3187.5
$ws_LTW.Range("M68").Value = -21329
$ws_LTW.Range("N68").Value = -4685.5
$ws_LTW.Range("H71").Value = 8158.684
$ws_LTW.Range("I71").Value = 22078
$ws_LTW.Range("J71").Value = 3187.5
$ws_LTW.Range("K71").Value = 110390
$ws_LTW.Range("L71").Value = 15937.5
$ws_LTW.Range("M71").Value = -106646
$ws_LTW.Range("N71").Value = -23425.5
$ws_LTW.Range("H82").Value = 1260.6666
$ws_LTW.Range("I82").Value = 1765.3334
$ws_LTW.Range("J82").Value = 1008.3333
$ws_LTW.Range("K82").Value = 1765.3334
$ws_LTW.Range("L82").Value = 1008.3333
$ws_LTW.Range("M82").Value = -1404.3334
$ws_LTW.Range("N82").Value = -1730.3333
$ws_LTW.Range("H85").Value = 1260.6666
$ws_LTW.Range("I85").Value = 1765.3334
$ws_LTW.Range("J85").Value = 1008.3333
$ws_LTW.Range("K85").Value = 1765.3334
$ws_LTW.Range("L85").Value = 1008.3333
$ws_LTW.Range("M85").Value = -517.3334
$ws_LTW.Range("N85").Value = -3504.3333
$ws_LTW.Range("H126").Value = 2624.75
$ws_LTW.Range("I126").Value = 2500
$ws_LTW.Range("J126").Value = 2999
$ws_LTW.Range("K126").Value = 7500
$ws_LTW.Range("L126").Value = 8997
$ws_LTW.Range("M126").Value = -5030
$ws_LTW.Range("N126").Value = -13937
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_WVR.Range("H62").Value = 2309.625
$ws_WVR.Range("J62").Value = 1939.5714
$ws_WVR.Range("L62").Value = 1939.5714
$ws_WVR.Range("N62").Value = -3187.5714
$ws_WVR.Range("H65").Value = 2309.625
$ws_WVR.Range("J65").Value = 1939.5714
$ws_WVR.Range("L65").Value = 9697.857
$ws_WVR.Range("N65").Value = -15937.857
$ws_WVR.Range("H81").Value = 931.6818
$ws_WVR.Range("I81").Value = 972.25
$ws_WVR.Range("J81").Value = 823.5
$ws_WVR.Range("K81").Value = 1944.5
$ws_WVR.Range("L81").Value = 1647
$ws_WVR.Range("M81").Value = -883.5
$ws_WVR.Range("N81").Value = -3769
$ws_WVR.Range("H84").Value = 931.6818
$ws_WVR.Range("I84").Value = 972.25
$ws_WVR.Range("J84").Value = 823.5
$ws_WVR.Range("K84").Value = 9722.5
$ws_WVR.Range("L84").Value = 8235
$ws_WVR.Range("M84").Value = -4418.5
$ws_WVR.Range("N84").Value = -18843
$ws_WVR.Range("H122").Value = 71187.05
$ws_WVR.Range("I122").Value = 17612.545
$ws_WVR.Range("J122").Value = 136667
$ws_WVR.Range("K122").Value = 52837.63499999999
$ws_WVR.Range("L122").Value = 410001
$ws_WVR.Range("M122").Value = -50387.63499999999
$ws_WVR.Range("N122").Value = -414901
$ws_WVR.Range("H132").Value = 2995.65
$ws_WVR.Range("I132").Value = 3665.125
$ws_WVR.Range("J132").Value = 1991.4375
$ws_WVR.Range("K132").Value = 10995.375
$ws_WVR.Range("L132").Value = 5974.3125
$ws_WVR.Range("M132").Value = -8465.375
$ws_WVR.Range("N132").Value = -11034.3125
